$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 45116.91341651211
$ws.Range("F2").Value = 45116.91341651322

$ws.Range("E3").Value = 45116.91341999901
$ws.Range("F3").Value = 45116.91342195809
$ws.Range("G3").Value = "0d, 0hr, 0min, 0.169sec "

$ws.Range("E4").Value = 45116.91342219352
$ws.Range("F4").Value = 45116.91342400664
$ws.Range("G4").Value = "0d, 0hr, 0min, 0.156sec "

$ws.Range("E5").Value = 45116.91342430653
$ws.Range("F5").Value = 45116.91342512117
$ws.Range("G5").Value = "0d, 0hr, 0min, 0.07sec "

$ws.Range("E6").Value = 45116.91342541348
$ws.Range("F6").Value = 45116.91342563316
$ws.Range("G6").Value = "0d, 0hr, 0min, 0.018sec "

$ws.Range("E7").Value = 45116.91342593425
$ws.Range("F7").Value = 45116.91342863093
$ws.Range("G7").Value = "0d, 0hr, 0min, 0.232sec "

$ws.Range("E8").Value = 45116.91342901658
$ws.Range("F8").Value = 45116.91343417042
$ws.Range("G8").Value = "0d, 0hr, 0min, 0.445sec "

$ws.Range("E9").Value = 45116.91343447495
$ws.Range("F9").Value = 45116.9134392728
$ws.Range("G9").Value = "0d, 0hr, 0min, 0.414sec "

$ws.Range("E10").Value = 45116.91343968582
$ws.Range("F10").Value = 45116.9134458843
$ws.Range("G10").Value = "0d, 0hr, 0min, 0.535sec "
